$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (adopts the varying values formerly on row 22)
$ws.Range("D2").Value = 45037
$ws.Range("N2").Value = 16000
$ws.Range("O2").Value = 16000
$ws.Range("P2").Value = 16000
$ws.Range("S2").Value = 889

# Row 3 (adopts the varying values formerly on row 8)
$ws.Range("D3").Value = 45001
$ws.Range("M3").Value = 60
$ws.Range("N3").Value = 17000
$ws.Range("O3").Value = 18000
$ws.Range("P3").Value = 17500
$ws.Range("S3").Value = 972

# Row 5 (adopts the varying values formerly on row 7)
$ws.Range("D5").Value = 45020
$ws.Range("N5").Value = 15000
$ws.Range("O5").Value = 15000
$ws.Range("P5").Value = 15000
$ws.Range("Q5").Value = '$/caja 16 kilos'
$ws.Range("R5").Value = 'Provincia de Los Andes'
$ws.Range("S5").Value = 938
$ws.Range("T5").Value = 16

# Row 6 (adopts the varying values formerly on row 16)
$ws.Range("D6").Value = 45044
$ws.Range("M6").Value = 60
$ws.Range("N6").Value = 15000
$ws.Range("O6").Value = 15000
$ws.Range("P6").Value = 15000
$ws.Range("S6").Value = 833

# Row 7 (adopts the varying values formerly on row 12)
$ws.Range("D7").Value = 45002
$ws.Range("M7").Value = 30
$ws.Range("N7").Value = 18000
$ws.Range("O7").Value = 18000
$ws.Range("P7").Value = 18000
$ws.Range("Q7").Value = '$/caja 18 kilos'
$ws.Range("R7").Value = 'Región Metropolitana'
$ws.Range("S7").Value = 1000
$ws.Range("T7").Value = 18

# Row 8 (adopts the varying values formerly on row 17)
$ws.Range("D8").Value = 44999

# Row 9 (adopts the varying values formerly on row 14)
$ws.Range("D9").Value = 45049
$ws.Range("M9").Value = 80
$ws.Range("N9").Value = 15000
$ws.Range("O9").Value = 15000
$ws.Range("P9").Value = 15000
$ws.Range("S9").Value = 833

# Row 10 (adopts the varying values formerly on row 20)
$ws.Range("D10").Value = 45021
$ws.Range("M10").Value = 60
$ws.Range("N10").Value = 15000
$ws.Range("O10").Value = 16000
$ws.Range("P10").Value = 15500
$ws.Range("R10").Value = 'Provincia de Los Andes'
$ws.Range("S10").Value = 861

# Row 11 (adopts the varying values formerly on row 18)
$ws.Range("D11").Value = 45062
$ws.Range("M11").Value = 90
$ws.Range("N11").Value = 13000
$ws.Range("O11").Value = 14000
$ws.Range("P11").Value = 13444
$ws.Range("S11").Value = 747

# Row 12 (adopts the varying values formerly on row 3)
$ws.Range("D12").Value = 45050
$ws.Range("M12").Value = 40
$ws.Range("N12").Value = 14000
$ws.Range("O12").Value = 14000
$ws.Range("P12").Value = 14000
$ws.Range("S12").Value = 778

# Row 13 (adopts the varying values formerly on row 5)
$ws.Range("D13").Value = 45096
$ws.Range("M13").Value = 50
$ws.Range("N13").Value = 23000
$ws.Range("O13").Value = 23000
$ws.Range("P13").Value = 23000
$ws.Range("S13").Value = 1278

# Row 14 (adopts the varying values formerly on row 19)
$ws.Range("D14").Value = 45030
$ws.Range("M14").Value = 40
$ws.Range("N14").Value = 18000
$ws.Range("O14").Value = 18000
$ws.Range("P14").Value = 18000
$ws.Range("S14").Value = 1000

# Row 15 (adopts the varying values formerly on row 13)
$ws.Range("D15").Value = 45033
$ws.Range("M15").Value = 60
$ws.Range("N15").Value = 15000
$ws.Range("O15").Value = 16000
$ws.Range("P15").Value = 15500
$ws.Range("S15").Value = 861

# Row 16 (adopts the varying values formerly on row 15)
$ws.Range("D16").Value = 45099
$ws.Range("M16").Value = 40
$ws.Range("N16").Value = 22000
$ws.Range("O16").Value = 22000
$ws.Range("P16").Value = 22000
$ws.Range("S16").Value = 1222

# Row 17 (adopts the varying values formerly on row 2)
$ws.Range("D17").Value = 45041
$ws.Range("N17").Value = 15000
$ws.Range("O17").Value = 15000
$ws.Range("P17").Value = 15000
$ws.Range("S17").Value = 833

# Row 18 (adopts the varying values formerly on row 11)
$ws.Range("D18").Value = 45028
$ws.Range("M18").Value = 50
$ws.Range("N18").Value = 18000
$ws.Range("O18").Value = 18000
$ws.Range("P18").Value = 18000
$ws.Range("S18").Value = 1000

# Row 19 (adopts the varying values formerly on row 10)
$ws.Range("D19").Value = 45091
$ws.Range("M19").Value = 50
$ws.Range("N19").Value = 22000
$ws.Range("O19").Value = 22000
$ws.Range("P19").Value = 22000
$ws.Range("S19").Value = 1222

# Row 20 (adopts the varying values formerly on row 6)
$ws.Range("D20").Value = 45014
$ws.Range("M20").Value = 30
$ws.Range("N20").Value = 18000
$ws.Range("O20").Value = 18000
$ws.Range("P20").Value = 18000
$ws.Range("R20").Value = 'Región Metropolitana'
$ws.Range("S20").Value = 1000

# Row 21 (adopts the varying values formerly on row 9)
$ws.Range("D21").Value = 45089
$ws.Range("N21").Value = 22000
$ws.Range("O21").Value = 23000
$ws.Range("P21").Value = 22500
$ws.Range("S21").Value = 1250

# Row 22 (adopts the varying values formerly on row 21)
$ws.Range("D22").Value = 45043
$ws.Range("N22").Value = 15000
$ws.Range("O22").Value = 15000
$ws.Range("P22").Value = 15000
$ws.Range("S22").Value = 833
